$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.413.31'
$ws.Range('E2').Value = '  -0.57%  '

$ws.Range('D3').Value = '2.610.97'
$ws.Range('E3').Value = '  -0.49%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').Value = '''591.98'
$ws.Range('E5').Value = '  -2.02%  '

$ws.Range('D6').Value = '''150.72'
$ws.Range('E6').Value = '  -2.68%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  -0.62%  '

$ws.Range('D9').Value = '2.609.80'
$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('E11').Value = '  -0.05%  '

$ws.Range('E13').Value = '  -3.48%  '

$ws.Range('D14').Value = '''27.30'
$ws.Range('E14').Value = '  -3.03%  '

$ws.Range('D15').Value = '3.081.62'
$ws.Range('E15').Value = '  -0.67%  '

$ws.Range('E16').Value = '  -2.67%  '

$ws.Range('D17').Value = '67.306.02'
$ws.Range('E17').Value = '  -0.59%  '

$ws.Range('D18').Value = '2.610.01'
$ws.Range('E18').Value = '  -0.37%  '

$ws.Range('D19').Value = '''371.77'
$ws.Range('E19').Value = '  +1.84%  '

$ws.Range('E20').Value = '  -2.73%  '

$ws.Range('E21').Value = '  -3.55%  '

$ws.Range('D22').Value = '''4.28'
$ws.Range('E22').Value = '  -0.67%  '

$ws.Range('E23').Value = '  -4.80%  '

$ws.Range('E24').Value = '  -4.01%  '

$ws.Range('D25').Value = '''73.68'
$ws.Range('E25').Value = '  +4.89%  '

$ws.Range('D27').Value = '''9.91'
$ws.Range('E27').Value = '  -2.40%  '

$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = '''580.30'
$ws.Range('E29').Value = '  -1.21%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''0.999'
$ws.Range('E30').Value = '  -0.50%  '

$ws.Range('E31').Value = '  -6.60%  '

$ws.Range('E32').Value = '  -5.95%  '

$ws.Range('E33').Value = '  -3.74%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('E36').Value = '  -4.19%  '

$ws.Range('E37').Value = '  -2.95%  '

$ws.Range('D38').Value = '''157.90'
$ws.Range('E38').Value = '  +1.44%  '

$ws.Range('D39').Value = '''19.06'
$ws.Range('E39').Value = '  -2.17%  '

$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.87'
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '''0.366'
$ws.Range('E41').Value = '  -1.95%  '

$ws.Range('E42').Value = '  -3.95%  '

$ws.Range('E43').Value = '  -3.85%  '

$ws.Range('E44').Value = '  +4.25%  '

$ws.Range('E45').Value = '  +0.05%  '

$ws.Range('D46').Value = '''153.43'
$ws.Range('E46').Value = '  -2.57%  '

$ws.Range('D47').Value = '0.0₆0282'
$ws.Range('E47').Value = '  -2.25%  '

$ws.Range('E48').Value = '  -1.59%  '

$ws.Range('D49').Value = '''0.0779'
$ws.Range('E49').Value = '  -1.53%  '

$ws.Range('E50').Value = '  -6.19%  '

$ws.Range('D51').Value = '''21.37'
$ws.Range('E51').Value = '  +1.07%  '
